# Update "想去人数" (interest count) values across the four worksheets to
# reflect a newer data pull (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param(
        [string]$SheetName,
        [hashtable]$RowToValue
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowToValue.Keys) {
        $ws.Range("F$row").Value = $RowToValue[$row]
    }
}

# Sheet: 展览 (Exhibitions)
Set-FValues "展览" @{
    3  = 1083
    8  = 13163
    9  = 13163
    10 = 2272
    13 = 54026
    15 = 312
    16 = 304
    19 = 357
    20 = 2972
    21 = 860
    22 = 5133
    23 = 1249
    24 = 935
    28 = 373
    29 = 1196
    30 = 82
    32 = 146
    38 = 4715
    39 = 33
    40 = 4738
    41 = 8706
    44 = 117
    46 = 404
    47 = 98
    49 = 4160
    50 = 178
}

# Sheet: 演出 (Performances)
Set-FValues "演出" @{
    4  = 89
    12 = 1115
}

# Sheet: 本地生活 (Local Life)
Set-FValues "本地生活" @{
    2 = 780
    3 = 543
    5 = 31
}

# Sheet: 全部类型 (All Types)
Set-FValues "全部类型" @{
    2  = 780
    3  = 543
    5  = 1083
    9  = 13163
    10 = 2272
    12 = 304
    15 = 2972
    16 = 860
    17 = 89
    18 = 1249
    19 = 31
    20 = 935
    27 = 82
    28 = 146
    32 = 4715
    33 = 33
    34 = 4738
    38 = 117
    43 = 98
    45 = 4160
}
